{"js": "// Capitalize the \"S\" in \"Yuccasoft\" -> \"YuccaSoft\" within the\n// \"Junior Programmer\" heading line.\nconst searchResults = context.document.body.search(\"Yuccasoft\", { matchCase: true });\nsearchResults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < searchResults.items.length; i++) {\n  searchResults.items[i].insertText(\"YuccaSoft\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Fix capitalization: \"Yuccasoft\" -> \"YuccaSoft\" (capital S)\n# in the \"Junior Programmer\" experience heading.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"Yuccasoft\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"YuccaSoft\"\n\n# wdFindContinue=1, wdReplaceAll=2\n$find.Execute(\"Yuccasoft\", $true, $false, $false, $false, $false, $true, 1, $false, \"YuccaSoft\", 2)\n"}
